$d = $word.ActiveDocument

# Locate the "Summary" section's paragraph. Nearly identical wording also
# appears later in the "Background" section, so find the paragraph by its
# distinguishing opening text rather than relying on a fixed index, then
# scope the Find/Replace to that paragraph only.
$summaryPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Polyploidisation is thought to be fundamental in the diversification of plant species. Both mutation")) {
        $summaryPara = $p
        break
    }
}

$rng = $summaryPara.Range
$rng.Find.Execute("locations but disentangling", $true, $false, $false, $false, $false, $true, 1, $false, "locations, but  disentangling", 2)

# Append a trailing space at the very end of the (now longer) paragraph.
$rng2 = $summaryPara.Range
$insertionPoint = $d.Range($rng2.End - 1, $rng2.End - 1)
$insertionPoint.InsertAfter(" ")
